# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values for the last (ca3748de-...) row on both the zh-cn and de-de report
# sheets, reflecting a newer report generation run ("Generate Report for
# Handback").

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-27 08:28:02"
$wsZhCn.Range("G5").Value = "2016-01-27 08:28:59"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-27 08:28:17"
$wsDeDe.Range("G5").Value = "2016-01-27 08:29:21"
